# Auto-generated Excel COM-interop script
# Applies cached-value refresh for Siren_Profits leve-profit columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets, per the scheduled-runner diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 42
$ws.Range("H42").Value = 344.23077
$ws.Range("I42").Value = 71.333336
$ws.Range("J42").Value = 426.1
$ws.Range("K42").Value = 214.000008
$ws.Range("L42").Value = 1278.3
$ws.Range("M42").Value = 15.99999199999999
$ws.Range("N42").Value = -1738.3
# Row 106
$ws.Range("H106").Value = 4325.636
$ws.Range("I106").Value = 4944.6
$ws.Range("J106").Value = 2999.2856
$ws.Range("K106").Value = 4944.6
$ws.Range("L106").Value = 2999.2856
$ws.Range("M106").Value = -4313.6
# Row 140
$ws.Range("H140").Value = 55069.43
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 55069.43
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 55069.43
$ws.Range("N140").Value = -65429.43
# Row 141
$ws.Range("H141").Value = 7066.077
$ws.Range("I141").Value = 6485.9
$ws.Range("J141").Value = 9000
$ws.Range("K141").Value = 19457.7
$ws.Range("L141").Value = 27000
$ws.Range("M141").Value = -14277.7
$ws.Range("N141").Value = -37360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 129763750
$ws.Range("I28").Value = 18333
$ws.Range("J28").Value = 519000000
$ws.Range("K28").Value = 18333
$ws.Range("L28").Value = 519000000
$ws.Range("M28").Value = -18141
# Row 92
$ws.Range("H92").Value = 366686340
$ws.Range("I92").Value = 59000
$ws.Range("J92").Value = 550000000
$ws.Range("K92").Value = 59000
$ws.Range("L92").Value = 550000000
$ws.Range("M92").Value = -56504
$ws.Range("N92").Value = -550004992
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 94
$ws.Range("H94").Value = 330000000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 330000000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 330000000
$ws.Range("N94").Value = -330001802
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 97
$ws.Range("H97").Value = 9897.134
$ws.Range("I97").Value = 10121.929
$ws.Range("J97").Value = 6750
$ws.Range("K97").Value = 10121.929
$ws.Range("L97").Value = 6750
$ws.Range("M97").Value = -9625.929
# Row 99
$ws.Range("H99").Value = 129763750
$ws.Range("I99").Value = 18333
$ws.Range("J99").Value = 519000000
$ws.Range("K99").Value = 18333
$ws.Range("L99").Value = 519000000
$ws.Range("M99").Value = -15338

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4645.148
$ws.Range("I86").Value = 5878.8
$ws.Range("J86").Value = 3103.0833
$ws.Range("K86").Value = 5878.8
$ws.Range("L86").Value = 3103.0833
$ws.Range("M86").Value = -4755.8
$ws.Range("N86").Value = -5349.0833
# Row 89
$ws.Range("H89").Value = 4645.148
$ws.Range("I89").Value = 5878.8
$ws.Range("J89").Value = 3103.0833
$ws.Range("K89").Value = 29394
$ws.Range("L89").Value = 15515.4165
$ws.Range("M89").Value = -23778
$ws.Range("N89").Value = -26747.4165
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 94
$ws.Range("H94").Value = 4984.125
$ws.Range("I94").Value = 4895.5
$ws.Range("J94").Value = 5250
$ws.Range("K94").Value = 4895.5
$ws.Range("L94").Value = 5250
$ws.Range("M94").Value = -4444.5
$ws.Range("N94").Value = -6152
# Row 95
$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
# Row 96
$ws.Range("H96").Value = 19616.334
$ws.Range("I96").Value = 19616.334
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 19616.334
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -16870.334
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5099.923
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5099.923
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5099.923
$ws.Range("N31").Value = -5689.923
$ws.Range("M31").ClearContents()
# Row 34
$ws.Range("H34").Value = 5099.923
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5099.923
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5099.923
$ws.Range("N34").Value = -5503.923
$ws.Range("M34").ClearContents()
# Row 51
$ws.Range("H51").Value = 35000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 35000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 35000
$ws.Range("N51").Value = -36472
# Row 61
$ws.Range("H61").Value = 35000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 35000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 35000
$ws.Range("N61").Value = -35696
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 45204050
$ws.Range("I4").Value = 35700884
$ws.Range("J4").Value = 140235660
$ws.Range("K4").Value = 107102652
$ws.Range("L4").Value = 420706980
$ws.Range("M4").Value = -107102540
# Row 11
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 600
$ws.Range("N11").Value = -880
# Row 34
$ws.Range("H34").Value = 1138801.1
$ws.Range("I34").Value = 2780212
$ws.Range("J34").Value = 2439.7693
$ws.Range("K34").Value = 8340636
$ws.Range("L34").Value = 7319.3079
$ws.Range("M34").Value = -8340552
$ws.Range("N34").Value = -7487.3079
# Row 46
$ws.Range("H46").Value = 1566.3636
$ws.Range("I46").Value = 530.3333
$ws.Range("J46").Value = 2809.6
$ws.Range("K46").Value = 1590.9999
$ws.Range("L46").Value = 8428.799999999999
$ws.Range("M46").Value = -1499.9999
$ws.Range("N46").Value = -8610.799999999999
# Row 55
$ws.Range("H55").Value = 8880.223
$ws.Range("I55").Value = 1129.3334
$ws.Range("J55").Value = 9849.083000000001
$ws.Range("K55").Value = 3388.0002
$ws.Range("L55").Value = 29547.249
$ws.Range("M55").Value = -3211.0002
$ws.Range("N55").Value = -29901.249
# Row 61
$ws.Range("H61").Value = 2664.8333
$ws.Range("I61").Value = 2664.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7994.499899999999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7779.499899999999
# Row 68
$ws.Range("H68").Value = 41678484
$ws.Range("I68").Value = 1698.5
$ws.Range("J68").Value = 62516876
$ws.Range("K68").Value = 5095.5
$ws.Range("L68").Value = 187550628
$ws.Range("M68").Value = -4284.5
$ws.Range("N68").Value = -187552250
# Row 71
$ws.Range("H71").Value = 41678484
$ws.Range("I71").Value = 1698.5
$ws.Range("J71").Value = 62516876
$ws.Range("K71").Value = 15286.5
$ws.Range("L71").Value = 562651884
$ws.Range("M71").Value = -11230.5
$ws.Range("N71").Value = -562659996
# Row 80
$ws.Range("H80").Value = 244699.8
$ws.Range("I80").Value = 9999.5
$ws.Range("J80").Value = 401166.66
$ws.Range("K80").Value = 29998.5
$ws.Range("L80").Value = 1203499.98
$ws.Range("M80").Value = -29062.5
$ws.Range("N80").Value = -1205371.98
# Row 83
$ws.Range("H83").Value = 244699.8
$ws.Range("I83").Value = 9999.5
$ws.Range("J83").Value = 401166.66
$ws.Range("K83").Value = 89995.5
$ws.Range("L83").Value = 3610499.94
$ws.Range("M83").Value = -85315.5
$ws.Range("N83").Value = -3619859.94
# Row 92
$ws.Range("H92").Value = 10150
$ws.Range("I92").Value = 99.75
$ws.Range("J92").Value = 16850.166
$ws.Range("K92").Value = 299.25
$ws.Range("L92").Value = 50550.49800000001
$ws.Range("M92").Value = 948.75
$ws.Range("N92").Value = -53046.49800000001
# Row 97
$ws.Range("H97").Value = 24703.191
$ws.Range("I97").Value = 34886.11
$ws.Range("J97").Value = 1791.625
$ws.Range("K97").Value = 104658.33
$ws.Range("L97").Value = 5374.875
$ws.Range("M97").Value = -104162.33
$ws.Range("N97").Value = -6366.875
# Row 118
$ws.Range("H118").Value = 3136.3333
$ws.Range("I118").Value = 954.5
$ws.Range("J118").Value = 7500
$ws.Range("K118").Value = 2863.5
$ws.Range("L118").Value = 22500
$ws.Range("M118").Value = -1620.5
$ws.Range("N118").Value = -24986
# Row 129
$ws.Range("H129").Value = 3185.1428
$ws.Range("I129").Value = 721.5
$ws.Range("J129").Value = 5032.875
$ws.Range("K129").Value = 2164.5
$ws.Range("L129").Value = 15098.625
$ws.Range("M129").Value = 2835.5
$ws.Range("N129").Value = -25098.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2392.0476
$ws.Range("I22").Value = 3048.4546
$ws.Range("J22").Value = 1670
$ws.Range("K22").Value = 3048.4546
$ws.Range("L22").Value = 1670
$ws.Range("M22").Value = -2753.4546
$ws.Range("N22").Value = -2260
# Row 27
$ws.Range("H27").Value = 2392.0476
$ws.Range("I27").Value = 3048.4546
$ws.Range("J27").Value = 1670
$ws.Range("K27").Value = 3048.4546
$ws.Range("L27").Value = 1670
$ws.Range("M27").Value = -2941.4546
$ws.Range("N27").Value = -1884

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 2749.7144
$ws.Range("I113").Value = 982.3077
$ws.Range("J113").Value = 5621.75
$ws.Range("K113").Value = 2946.9231
$ws.Range("L113").Value = 16865.25
$ws.Range("M113").Value = -776.9231
$ws.Range("N113").Value = -21205.25
# Row 122
$ws.Range("H122").Value = 8457.5
$ws.Range("I122").Value = 4050.1428
$ws.Range("J122").Value = 10830.692
$ws.Range("K122").Value = 12150.4284
$ws.Range("L122").Value = 32492.076
$ws.Range("M122").Value = -9700.428400000001
# Row 123
$ws.Range("H123").Value = 29515
$ws.Range("I123").Value = 29515
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 29515
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -24615

Write-Host "Applied 271 cell updates and 10 cell clears across 7 sheets."
